$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New average row right under the per-instance data (|S*|/n average)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Summary block: averages / worst-case ratios
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Formatting for the new summary labels/values: bold, 12pt, vertically centered
$summary = $ws.Range("A14:B17")
$summary.Font.Bold = $true
$summary.Font.Size = 12
$summary.VerticalAlignment = -4108

# Select the summary block as the final selection, matching the saved view state
$ws.Range("A14:B17").Select() | Out-Null

# Page setup metadata present in the resaved workbook
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
